$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.547.05"
$ws.Range("E2").Value = "  -0.28%  "
$ws.Range("D3").Value = "1.834.37"
$ws.Range("E3").Value = "  -0.52%  "
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").Value = "'312.05"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -0.29%  "
$ws.Range("E6").Value = "  -0.03%  "
$ws.Range("D7").Value = "'0.4276"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  -0.05%  "
$ws.Range("D8").Value = "'0.3655"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  +0.59%  "
$ws.Range("D9").Value = "'0.07271"
$ws.Range("D9").ClearFormats()
$ws.Range("D10").Value = "'0.8641"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  -1.81%  "
$ws.Range("D11").Value = "'20.66"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  -0.01%  "
$ws.Range("B12").Value = "WrappedEther"
$ws.Range("C12").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D12").Value = "1.760.06"
$ws.Range("E12").Value = "  -4.75%  "
$ws.Range("B13").Value = "Polkadot"
$ws.Range("C13").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D13").Value = "'5.467"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  +2.30%  "
$ws.Range("D14").Value = "'6.520"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  -0.18%  "
$ws.Range("D15").Value = "'0.06969"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  +0.23%  "
$ws.Range("D16").Value = "'1.003"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  +0.01%  "
$ws.Range("D17").Value = "'80.66"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  +1.40%  "
$ws.Range("D18").Value = "'0.000008917"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  -0.60%  "
$ws.Range("E19").Value = "  -0.12%  "
$ws.Range("D20").Value = "'15.40"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  +0.26%  "
$ws.Range("D21").Value = "27.309.05"
$ws.Range("E21").Value = "  -1.17%  "
$ws.Range("D22").Value = "'5.152"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  +3.25%  "
$ws.Range("E23").Value = "  +5.59%  "
$ws.Range("D24").Value = "1.984.93"
$ws.Range("E24").Value = "  -4.26%  "
$ws.Range("D25").Value = "'1.992"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  +0.14%  "
$ws.Range("D26").Value = "'154.97"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  -0.31%  "
$ws.Range("E27").Value = "  +1.94%  "
$ws.Range("D28").Value = "'5.152"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  -1.69%  "
$ws.Range("D29").Value = "'114.23"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  -4.60%  "
$ws.Range("D30").Value = "'1.817"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  -3.58%  "
$ws.Range("D31").Value = "'0.08853"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  -0.41%  "
$ws.Range("D32").Value = "'0.7482"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  -2.09%  "
$ws.Range("D33").Value = "'2.987"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  +1.12%  "
$ws.Range("D34").Value = "'4.543"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  +0.44%  "
$ws.Range("D35").Value = "'1.132"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  +0.07%  "
$ws.Range("E36").Value = "  +0.02%  "
$ws.Range("D37").Value = "'1.097"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  -0.81%  "
$ws.Range("E38").Value = "  -2.67%  "
$ws.Range("D39").Value = "'0.01936"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  -0.13%  "
$ws.Range("D40").Value = "'2.801"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  -0.58%  "
$ws.Range("D41").Value = "'0.5068"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  -0.26%  "
$ws.Range("E42").Value = "  -1.00%  "
$ws.Range("D43").Value = "'6.451"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  -2.15%  "
$ws.Range("D44").Value = "'8.331"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  -0.97%  "
$ws.Range("D45").Value = "'10.40"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  -0.18%  "
$ws.Range("D46").Value = "'105.43"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  -0.44%  "
$ws.Range("D47").Value = "'0.06472"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  -1.17%  "
$ws.Range("D48").Value = "'0.4690"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  +0.96%  "
$ws.Range("D49").Value = "'1.000"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  -0.04%  "
$ws.Range("E50").Value = "  -1.11%  "
$ws.Range("D51").Value = "'1.740"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  -0.52%  "
